# Generate Report for Handback
# The two tracked files (dfddac69... and 9f0acc76...) have both now been
# handed back (in sync with en-US). The report rows are regenerated:
#  - file rows are reordered (9f0acc76 now sorts/appears before dfddac69)
#  - status text changes from "Ready for handoff" / mixed to
#    "Handed back: in sync with en-US" for both files
#  - the handback timestamps are refreshed to the new handback time

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: update a cell's value, and if that cell carries a hyperlink,
# also update the hyperlink's displayed text (without touching the
# hyperlink's target address, its relationship id, or the cell style).
# ---------------------------------------------------------------------
function Set-CellAndHyperlinkText {
    param(
        $ws,
        [string]$cellAddr,
        [string]$newText
    )

    $ws.Range($cellAddr).Value2 = $newText

    $target = '$' + ($cellAddr -replace '^([A-Z]+)(\d+)$', '$1$$$2')
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $target) {
            $hl.TextToDisplay = $newText
        }
    }
}

# =======================================================================
# Sheet "Overview"
# =======================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

Set-CellAndHyperlinkText $wsOverview "A2" "9f0acc76-86c5-4abe-a47a-94a2628b08ba.md"
Set-CellAndHyperlinkText $wsOverview "A3" "dfddac69-0996-44ef-9b10-44303e86f223.md"

$wsOverview.Range("B3").Value2 = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value2 = "Handed back: in sync with en-US"

# =======================================================================
# Sheet "zh-cn"
# =======================================================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-CellAndHyperlinkText $wsZhCn "A2" "9f0acc76-86c5-4abe-a47a-94a2628b08ba.md"
Set-CellAndHyperlinkText $wsZhCn "E2" "9f0acc76-86c5-4abe-a47a-94a2628b08ba.md"
Set-CellAndHyperlinkText $wsZhCn "C2" "9f0acc76-86c5-4abe-a47a-94a2628b08ba.359142bc38900f7c04cccebb26aee01a967cf849.zh-cn.xlf"
Set-CellAndHyperlinkText $wsZhCn "F2" "9f0acc76-86c5-4abe-a47a-94a2628b08ba.359142bc38900f7c04cccebb26aee01a967cf849.zh-cn.xlf"

Set-CellAndHyperlinkText $wsZhCn "A3" "dfddac69-0996-44ef-9b10-44303e86f223.md"
Set-CellAndHyperlinkText $wsZhCn "E3" "dfddac69-0996-44ef-9b10-44303e86f223.md"
Set-CellAndHyperlinkText $wsZhCn "C3" "dfddac69-0996-44ef-9b10-44303e86f223.80679bb79595f5dbf74c41c01399de422e2a60ff.zh-cn.xlf"
Set-CellAndHyperlinkText $wsZhCn "F3" "dfddac69-0996-44ef-9b10-44303e86f223.80679bb79595f5dbf74c41c01399de422e2a60ff.zh-cn.xlf"

$wsZhCn.Range("B3").Value2 = "Handed back: in sync with en-US"
$wsZhCn.Range("G2").Value2 = "2016-02-29 04:37:13"
$wsZhCn.Range("G3").Value2 = "2016-02-29 04:37:13"

# =======================================================================
# Sheet "de-de"
# =======================================================================
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-CellAndHyperlinkText $wsDeDe "A2" "9f0acc76-86c5-4abe-a47a-94a2628b08ba.md"
Set-CellAndHyperlinkText $wsDeDe "E2" "9f0acc76-86c5-4abe-a47a-94a2628b08ba.md"
Set-CellAndHyperlinkText $wsDeDe "C2" "9f0acc76-86c5-4abe-a47a-94a2628b08ba.359142bc38900f7c04cccebb26aee01a967cf849.de-de.xlf"
Set-CellAndHyperlinkText $wsDeDe "F2" "9f0acc76-86c5-4abe-a47a-94a2628b08ba.359142bc38900f7c04cccebb26aee01a967cf849.de-de.xlf"

Set-CellAndHyperlinkText $wsDeDe "A3" "dfddac69-0996-44ef-9b10-44303e86f223.md"
Set-CellAndHyperlinkText $wsDeDe "E3" "dfddac69-0996-44ef-9b10-44303e86f223.md"
Set-CellAndHyperlinkText $wsDeDe "C3" "dfddac69-0996-44ef-9b10-44303e86f223.80679bb79595f5dbf74c41c01399de422e2a60ff.de-de.xlf"
Set-CellAndHyperlinkText $wsDeDe "F3" "dfddac69-0996-44ef-9b10-44303e86f223.80679bb79595f5dbf74c41c01399de422e2a60ff.de-de.xlf"

$wsDeDe.Range("B3").Value2 = "Handed back: in sync with en-US"
$wsDeDe.Range("G2").Value2 = "2016-02-29 04:37:35"
$wsDeDe.Range("G3").Value2 = "2016-02-29 04:37:35"
